# Removed Test Case Inter-Dependency
#
# - ProductLoanInput!B1 (productname) and ProductLoanOutput!B1 (verifyloanproduct)
#   no longer hard-code the shared "4350-Simple-Group-Loan-Product" name; they now
#   use a distinct "-1st" suffixed name so the two test sheets don't collide.
# - ProductLoanInput!B2 (shortname) switches from the numeric 4350 to a unique
#   text short name "435p".
# - Active sheet/selection moves from ProductLoanInput!B18 to ProductLoanOutput
#   (which becomes the selected tab), with ProductLoanInput's own selection
#   reset to B3.

$wb  = $excel.ActiveWorkbook
$wsIn  = $wb.Worksheets.Item("ProductLoanInput")
$wsOut = $wb.Worksheets.Item("ProductLoanOutput")

# New product name used by both sheets, and a new text shortname.
$wsIn.Range("B1").Value = "4350-Simple-Group-Loan-Product-1st"
$wsIn.Range("B2").Value = "435p"
$wsOut.Range("B1").Value = "4350-Simple-Group-Loan-Product-1st"

# Reset the input sheet's selection (it is no longer the active tab).
$wsIn.Range("B3").Select()

# Output sheet becomes the active/selected tab, selection stays at B1.
$wsOut.Activate()
$wsOut.Range("B1").Select()
